$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35: new journal entry - 14/04/2018, crud/ClientRepository/ORM work, 0.5h
$ws.Range("A35").Value = 43204
$ws.Range("A35").NumberFormat = "mm-dd-yy"
$ws.Range("B35").Value = "développement du crud et test sur ClientRepository et mise en place de la classe ORM"
$ws.Range("C35").Value = 0.5

# Row 36: new journal entry - 14/04/2018, ClientRepository derby tests (no hours yet)
$ws.Range("A36").Value = 43204
$ws.Range("B36").Value = "développement et test de ClientRepository pour derby"

# Reuse A35's date format/style for A36 instead of creating a duplicate style
$ws.Range("A35").Copy()
$ws.Range("A36").PasteSpecial(-4122)

# Move the active selection to where the user left off editing
$ws.Range("J31").Select()
